$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map (2024-02-10 cryptos refresh).
# Values are apostrophe-prefixed so Excel stores them as literal text
# (matching the original inlineStr cells) instead of coercing numeric-
# looking strings like "0.524" into real numbers; Style is reset to
# "Normal" afterwards so the quote-prefix flag does not linger in the
# saved cell format.
$updates = [ordered]@{
    'D2' = '47.882.80'
    'E2' = '  +1.29%  '
    'D3' = '2.501.31'
    'E3' = '  +0.21%  '
    'E4' = '  +0.03%  '
    'D5' = '323.52'
    'E5' = '  -0.17%  '
    'D6' = '109.22'
    'E6' = '  +1.19%  '
    'D7' = '0.524'
    'E7' = '  -0.41%  '
    'E8' = '  +0.05%  '
    'D9' = '0.552'
    'E9' = '  +1.84%  '
    'D10' = '40.44'
    'E10' = '  +5.53%  '
    'D11' = '0.0815'
    'E11' = '  -0.12%  '
    'E12' = '  +0.41%  '
    'D13' = '18.93'
    'E13' = '  +2.13%  '
    'D14' = '7.23'
    'E14' = '  +0.10%  '
    'D15' = '2.892.55'
    'E15' = '  +0.27%  '
    'D16' = '2.497.94'
    'E16' = '  -0.15%  '
    'D17' = '0.852'
    'E17' = '  -0.71%  '
    'D18' = '47.774.90'
    'E18' = '  +1.21%  '
    'D19' = '13.20'
    'E19' = '  +1.51%  '
    'D20' = '6.64'
    'E20' = '  -0.86%  '
    'E21' = '  +12.65%  '
    'D22' = '0.0₃0943'
    'E22' = '  +0.13%  '
    'D23' = '70.82'
    'E23' = '  -0.02%  '
    'D24' = '248.10'
    'E24' = '  -1.11%  '
    'D25' = '2.57'
    'E25' = '  -1.18%  '
    'D27' = '25.94'
    'E27' = '  -1.29%  '
    'D28' = '10.00'
    'E28' = '  -0.49%  '
    'D29' = '2.21'
    'E29' = '  -3.79%  '
    'D30' = '35.24'
    'E30' = '  +0.56%  '
    'E31' = '  +0.74%  '
    'D32' = '49.89'
    'E32' = '  +0.86%  '
    'D33' = '19.89'
    'E33' = '  +0.82%  '
    'D34' = '5.36'
    'E34' = '  -2.33%  '
    'D35' = '0.0792'
    'E35' = '  -0.21%  '
    'E36' = '  +0.12%  '
    'D37' = '1.96'
    'E37' = '  -1.27%  '
    'D38' = '4.67'
    'E38' = '  -1.08%  '
    'E39' = '  -0.81%  '
    'E40' = '  -0.33%  '
    'D41' = '22.30'
    'E41' = '  +5.48%  '
    'E42' = '  -1.08%  '
    'D43' = '119.29'
    'E43' = '  -2.57%  '
    'E44' = '  -0.24%  '
    'D45' = '2.002.32'
    'E45' = '  +1.66%  '
    'D46' = '3.06'
    'E46' = '  +1.30%  '
    'E47' = '  -3.39%  '
    'D48' = '1.81'
    'E48' = '  +0.59%  '
    'D49' = '9.04'
    'E49' = '  -0.16%  '
    'D50' = '5.13'
    'E50' = '  -2.97%  '
    'D51' = '56.95'
    'E51' = '  +2.91%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}

